# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# (and the derived "Latest HO Xliff Generate Date" on the Overview sheet)
# for the 9f27646a-ddc7-4349-93bb-ac079887f00e file, on every localized
# language sheet, reflecting the freshly generated handback xliffs.

$wb = $excel.ActiveWorkbook

# Overview sheet: row 3 is the 9f27646a...md file; column G is
# "Latest HO Xliff Generate Date".
$overview = $wb.Worksheets("Overview")
$overview.Range("G3").Value = "2016-08-24 10:49:40"

# zh-cn sheet: row 3 is the 9f27646a...md file; H = "Correspond Handoff
# Datetime", K = "Correspond Handback DateTime".
$zhcn = $wb.Worksheets("zh-cn")
$zhcn.Range("H3").Value = "2016-08-24 10:49:34"
$zhcn.Range("K3").Value = "2016-08-24 10:49:52"

# de-de sheet: row 3 is the 9f27646a...md file; H = "Correspond Handoff
# Datetime", K = "Correspond Handback DateTime".
$dede = $wb.Worksheets("de-de")
$dede.Range("H3").Value = "2016-08-24 10:49:40"
$dede.Range("K3").Value = "2016-08-24 10:49:59"
